# Updated cryptos list (prices + 1h volume %) refreshed by the scraper run.
# Rows 20/21 and 41/42 additionally swap places (ranking reshuffle), so the
# Coin name + Link columns are rewritten for those rows too.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.576.97'
$ws.Range("E2").Value = '  +0.13%  '
$ws.Range("D3").Value = '1.961.59'
$ws.Range("E3").Value = '  +2.33%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9979'
$ws.Range("E4").Value = '  -0.38%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.66'
$ws.Range("E5").Value = '  +0.97%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9984'
$ws.Range("E6").Value = '  -0.30%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4805'
$ws.Range("E7").Value = '  -0.33%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2917'
$ws.Range("E8").Value = '  +0.59%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06752'
$ws.Range("E9").Value = '  +0.77%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '108.72'
$ws.Range("E10").Value = '  -1.91%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.18'
$ws.Range("E11").Value = '  +1.42%  '
$ws.Range("D12").Value = '1.955.58'
$ws.Range("E12").Value = '  +1.91%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07704'
$ws.Range("E13").Value = '  +1.91%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.396'
$ws.Range("E14").Value = '  +2.25%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6885'
$ws.Range("E15").Value = '  +3.06%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '289.76'
$ws.Range("E16").Value = '  -2.70%  '
$ws.Range("D17").Value = '30.591.26'
$ws.Range("E17").Value = '  +0.16%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.13'
$ws.Range("E18").Value = '  +1.27%  '
$ws.Range("D19").Value = '2.219.03'
$ws.Range("E19").Value = '  +2.24%  '
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007659'
$ws.Range("E20").Value = '  +1.13%  '
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.619'
$ws.Range("E21").Value = '  +1.40%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9980'
$ws.Range("E22").Value = '  -0.26%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9979'
$ws.Range("E23").Value = '  -0.35%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.610'
$ws.Range("E24").Value = '  +3.07%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.864'
$ws.Range("E25").Value = '  +4.42%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '171.07'
$ws.Range("E26").Value = '  +3.53%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.92'
$ws.Range("E27").Value = '  -1.40%  '
$ws.Range("E28").Value = '  +3.23%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1065'
$ws.Range("E29").Value = '  +0.24%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.438'
$ws.Range("E30").Value = '  +0.43%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.742'
$ws.Range("E31").Value = '  +16.88%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.443'
$ws.Range("E32").Value = '  +7.50%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05072'
$ws.Range("E33").Value = '  +1.42%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7661'
$ws.Range("E34").Value = '  +3.89%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.173'
$ws.Range("E35").Value = '  +3.43%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.726'
$ws.Range("E36").Value = '  +0.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02023'
$ws.Range("E37").Value = '  +0.54%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.711'
$ws.Range("E38").Value = '  +1.08%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.475'
$ws.Range("E39").Value = '  +11.10%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.126'
$ws.Range("E40").Value = '  +5.72%  '
$ws.Range("B41").Value = 'Quant'
$ws.Range("C41").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '109.38'
$ws.Range("E41").Value = '  -1.09%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8787'
$ws.Range("E42").Value = '  +1.60%  '
$ws.Range("E43").Value = '  +0.46%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '69.69'
$ws.Range("E44").Value = '  -0.96%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9980'
$ws.Range("E45").Value = '  -0.27%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.458'
$ws.Range("E46").Value = '  +3.64%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1264'
$ws.Range("E47").Value = '  +2.94%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.342'
$ws.Range("E48").Value = '  +1.01%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '35.75'
$ws.Range("E49").Value = '  +2.77%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '924.08'
$ws.Range("E50").Value = '  +6.57%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '47.06'
$ws.Range("E51").Value = '  -3.35%  '
